$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'43.535.97"
$ws.Range("E2").Value = "  -3.59%  "

$ws.Range("D3").Value = "'2.234.19"
$ws.Range("E3").Value = "  -5.21%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'319.10"
$ws.Range("E5").Value = "  +3.37%  "

$ws.Range("E6").Value = "  -7.42%  "

$ws.Range("D7").Value = "'0.583"
$ws.Range("E7").Value = "  -7.13%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "'0.569"
$ws.Range("E9").Value = "  -7.59%  "

$ws.Range("D10").Value = "'37.04"
$ws.Range("E10").Value = "  -10.11%  "

$ws.Range("D11").Value = "'54.45"
$ws.Range("E11").Value = "  -1.65%  "

$ws.Range("D12").Value = "'0.0829"
$ws.Range("E12").Value = "  -9.56%  "

$ws.Range("D13").Value = "'7.76"
$ws.Range("E13").Value = "  -8.05%  "

$ws.Range("E14").Value = "  -2.78%  "

$ws.Range("D15").Value = "'0.867"
$ws.Range("E15").Value = "  -11.77%  "

$ws.Range("D16").Value = "'2.573.27"
$ws.Range("E16").Value = "  -5.14%  "

$ws.Range("D17").Value = "'14.18"
$ws.Range("E17").Value = "  -7.40%  "

$ws.Range("D18").Value = "'2.227.86"
$ws.Range("E18").Value = "  -5.04%  "

$ws.Range("D19").Value = "'43.258.85"
$ws.Range("E19").Value = "  -4.00%  "

$ws.Range("D20").Value = "'14.42"
$ws.Range("E20").Value = "  +7.54%  "

$ws.Range("D21").Value = "'0.0₃0982"
$ws.Range("E21").Value = "  -7.83%  "

$ws.Range("D22").Value = "'6.56"
$ws.Range("E22").Value = "  -9.83%  "

$ws.Range("E23").Value = "  -9.96%  "

$ws.Range("D24").Value = "'3.22"
$ws.Range("E24").Value = "  -5.66%  "

$ws.Range("D25").Value = "'237.85"
$ws.Range("E25").Value = "  -8.15%  "

$ws.Range("D26").Value = "'2.16"
$ws.Range("E26").Value = "  -6.47%  "

$ws.Range("E27").Value = "  +0.17%  "

$ws.Range("D28").Value = "'10.27"
$ws.Range("E28").Value = "  -7.21%  "

$ws.Range("E29").Value = "  -6.77%  "

$ws.Range("D30").Value = "'6.50"
$ws.Range("E30").Value = "  -11.73%  "

$ws.Range("D31").Value = "'0.0897"
$ws.Range("E31").Value = "  -6.59%  "

$ws.Range("D32").Value = "'20.75"
$ws.Range("E32").Value = "  -7.36%  "

$ws.Range("D33").Value = "'34.43"
$ws.Range("E33").Value = "  -9.68%  "

$ws.Range("D34").Value = "'158.64"
$ws.Range("E34").Value = "  -6.99%  "

$ws.Range("D35").Value = "'2.78"
$ws.Range("E35").Value = "  -5.54%  "

$ws.Range("D36").Value = "'3.37"
$ws.Range("E36").Value = "  +14.78%  "

$ws.Range("E37").Value = "  -6.14%  "

$ws.Range("B38").Value = "ARBITRUM"
$ws.Range("C38").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D38").Value = "'1.92"
$ws.Range("E38").Value = "  +11.46%  "

$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'4.55"
$ws.Range("E39").Value = "  -6.29%  "

$ws.Range("D40").Value = "'0.105"
$ws.Range("E40").Value = "  -7.49%  "

$ws.Range("D41").Value = "'3.66"
$ws.Range("E41").Value = "  -6.33%  "

$ws.Range("E42").Value = "  -8.31%  "

$ws.Range("E43").Value = "  +0.03%  "

$ws.Range("D44").Value = "'1.826.18"
$ws.Range("E44").Value = "  +11.52%  "

$ws.Range("D45").Value = "'12.16"
$ws.Range("E45").Value = "  -3.97%  "

$ws.Range("D46").Value = "'89.33"
$ws.Range("E46").Value = "  -10.16%  "

$ws.Range("D47").Value = "'80.17"
$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'16.83"
$ws.Range("E48").Value = "  +69.59%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.209"
$ws.Range("E49").Value = "  -9.70%  "

$ws.Range("D50").Value = "'5.50"
$ws.Range("E50").Value = "  +0.10%  "

$ws.Range("D51").Value = "'61.53"
$ws.Range("E51").Value = "  -11.77%  "
